$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7 ("Experimental"): set the Value column (B7) to the text "true" ---
# A plain Value = "true" assignment gets auto-coerced by the engine into a
# boolean cell (t="b"). To store it as literal text (matching the target
# shared-string cell) we compute it as a text formula in a scratch cell,
# then paste just the resulting value into B7. This keeps B7's existing
# number format / style (s="2") untouched and avoids creating any stray
# cell styles.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="true"'
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$scratch.ClearContents()

# --- Row 8 ("Date"): update the recorded timestamp ---
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"

$excel.CutCopyMode = 0
